$d = $word.ActiveDocument

# The first paragraph of the document starts with a run containing just the
# letter "T" (the start of "TRUONG DAI HOC CAN THO"). The edit inserts a new
# run - containing only a tab character (a real <w:tab/> run, not a literal
# tab glyph inside a text run) - immediately before that existing run, using
# the same run formatting (bCs / color 000000 / lang vi-VN) already used by
# the paragraph's first run.

$target = $d.Paragraphs(1).Range
$insertionPoint = $target.Duplicate
$insertionPoint.SetRange($target.Start, $target.Start)

# Build a minimal single-part WordprocessingML package so the engine inserts
# an actual <w:tab/> element (InsertXML replaces the addressed - here
# zero-length - range with the supplied markup).
$tabRunPackage = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:bCs/>' +
                  '<w:color w:val="000000"/>' +
                  '<w:lang w:val="vi-VN"/>' +
                '</w:rPr>' +
                '<w:tab/>' +
              '</w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$null = $insertionPoint.InsertXML($tabRunPackage)
